# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5493
$ws1.Range("F4").Value = 636
$ws1.Range("F6").Value = 828
$ws1.Range("F7").Value = 30
$ws1.Range("F8").Value = 356
$ws1.Range("F10").Value = 1
$ws1.Range("F11").Value = 17

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5493
$ws4.Range("F4").Value = 636
$ws4.Range("F6").Value = 828
$ws4.Range("F7").Value = 30
$ws4.Range("F9").Value = 356
$ws4.Range("F11").Value = 1
$ws4.Range("F12").Value = 17

$wb.Save()
